# Add a "2022-Q4" sheet (fund-holding detail) right after the "总计"
# summary sheet, and update the "总计" table with a new leading row for
# the 2022-Q4 quarter (existing rows shift down).

$wb = $excel.ActiveWorkbook

# Helper: write a text-typed value (matches the inlineStr/text columns
# used for fund codes / 基金规模 / 股票总仓位 / 仓位占比 / 持有市值 throughout
# this workbook) by round-tripping through a text number-format, then
# clearing the format so no stray numFmt style sticks to the cell.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# ------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new row 2 for 2022-Q4
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Push existing data rows (2022-Q2, 2021-Q4, 2021-Q3) down by one row.
$summary.Rows.Item(2).Insert()

# The freshly inserted row picked up some inherited formatting from the
# header row above it - strip it back to the plain/no-style look used by
# every other data row in this table.
$summary.Range("B2:D2").ClearFormats()

# Column A keeps the bordered/bold "index" style used throughout the
# table (same style as A3/A4/A5) - grab it from the row below.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new 2022-Q4 summary row.
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 3
$summary.Range("D2").Value = 0

# Re-number the index column and refresh the values for the rows that
# shifted down (2022-Q2, 2021-Q4, 2021-Q3).
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2021-Q4"
$summary.Range("C4").Value = 4
$summary.Range("D4").Value = 0.97

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q3"
$summary.Range("C5").Value = 7
$summary.Range("D5").Value = 0.18

# ------------------------------------------------------------------
# 2. Insert the new "2022-Q4" worksheet right after "总计"
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $summary)
$newSheet.Name = "2022-Q4"

# Borrow the header-row (B1:H1) and index-column (A2) formatting from
# the "2022-Q2" sheet, which uses the same layout/style as the new tab.
# (Re-fetch the sheet reference now, AFTER the Add() above, so it points
# at the live object.)
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$q2Sheet.Range("A2").Copy()
$newSheet.Range("A2:A4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Headers
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2 - 007257 / 凯石沣混合A
$newSheet.Range("A2").Value = 0
Set-TextValue $newSheet.Range("B2") "007257"
$newSheet.Range("C2").Value = "凯石沣混合A"
Set-TextValue $newSheet.Range("D2") "0.08"
Set-TextValue $newSheet.Range("E2") "73.69"
Set-TextValue $newSheet.Range("F2") "2.55"
Set-TextValue $newSheet.Range("G2") "0.0020"
$newSheet.Range("H2").Value = 8

# Row 3 - 001797 / 国新国证新利灵活配置混合
$newSheet.Range("A3").Value = 1
Set-TextValue $newSheet.Range("B3") "001797"
$newSheet.Range("C3").Value = "国新国证新利灵活配置混合"
Set-TextValue $newSheet.Range("D3") "0.02"
Set-TextValue $newSheet.Range("E3") "81.37"
Set-TextValue $newSheet.Range("F3") "5.56"
Set-TextValue $newSheet.Range("G3") "0.0011"
$newSheet.Range("H3").Value = 3

# Row 4 - 007258 / 凯石沣混合C
$newSheet.Range("A4").Value = 2
Set-TextValue $newSheet.Range("B4") "007258"
$newSheet.Range("C4").Value = "凯石沣混合C"
Set-TextValue $newSheet.Range("D4") "0.04"
Set-TextValue $newSheet.Range("E4") "73.69"
Set-TextValue $newSheet.Range("F4") "2.55"
Set-TextValue $newSheet.Range("G4") "0.0010"
$newSheet.Range("H4").Value = 8

$newSheet.Range("A1").Select() | Out-Null

# Restore the original active sheet (总计 / summary tab) as the one
# shown when the workbook is opened, matching the untouched bookViews
# (activeTab=0) in the source file.
$summary = $wb.Worksheets.Item(1)
$summary.Activate()
$summary.Range("A1").Select() | Out-Null
